# Attendance sheet update:
#  - add new student "Cozariuc Vicentiu" with week 6 (column H) presence ticked
#  - tick week 6 (column H) presence for a batch of existing students
#  - re-sort the roster (B3:S42) alphabetically by first/last name (column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student goes into the first empty row (row 42) before the table is re-sorted
$ws.Range("B42").Value = "Cozariuc Vicentiu"
$ws.Range("H42").Value = $true

# Mark week 6 (sapt 6 / column H) attendance for the existing students who were present
$weekSixRows = @(6, 9, 14, 15, 16, 18, 20, 25, 31, 32, 37, 38, 41)
foreach ($r in $weekSixRows) {
    $ws.Cells.Item($r, 8).Value = $true
}

# Re-sort the whole roster (names + weekly marks) alphabetically by column B
$sortRange = $ws.Range("B3:S42")
$sortKey = $ws.Range("B42")
$sortRange.Sort($sortKey)
